# Update column F ("dSF") values on Sheet1 to match the repulled data.
# Per the commit "repull data, push all data, mean calculation" only the
# dSF column values for a subset of rows change; all other cells are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -4
    3  = -6
    6  = -1
    7  = 3
    10 = -4
    12 = 10
    13 = -13
    14 = 1
    17 = -2
    18 = -2
    19 = -4
    28 = 7
    30 = 3
    31 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
